$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 ("サキュバス") which pushes the current
# row 13 (ミニマム / minimum) down to row 14.
$ws.Rows.Item(13).Insert()

# New row 13: only column A is populated.
$ws.Cells.Item(13, 1).Value = "サキュバス"

# Fill column B for rows 5, 7, 8, 9, 10, 11, 12 with the same text as
# column A (these previously had no column-B / column-D counterpart).
$ws.Cells.Item(5, 2).Value = $ws.Cells.Item(5, 1).Text
$ws.Cells.Item(7, 2).Value = $ws.Cells.Item(7, 1).Text
$ws.Cells.Item(8, 2).Value = $ws.Cells.Item(8, 1).Text
$ws.Cells.Item(9, 2).Value = $ws.Cells.Item(9, 1).Text
$ws.Cells.Item(10, 2).Value = $ws.Cells.Item(10, 1).Text
$ws.Cells.Item(11, 2).Value = $ws.Cells.Item(11, 1).Text
$ws.Cells.Item(12, 2).Value = $ws.Cells.Item(12, 1).Text

# Column D was a (now redundant) mirror of column B; remove it entirely.
$ws.Columns.Item(4).Delete()
